$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-08-16"

# Update the August column header label (row 9, column A) to the new date
$ws.Range("A9").Value = "August (through 08-16)"

# Row 9 ("August (through 08-16)") - updated counts per year column
$ws.Range("B9").Value = 19
$ws.Range("C9").Value = 42
$ws.Range("D9").Value = 40
$ws.Range("E9").Value = 26
$ws.Range("F9").Value = 21
$ws.Range("G9").Value = 101
$ws.Range("H9").Value = 93
$ws.Range("I9").Value = 90

# Row 10 ("Total") - updated totals per year column
$ws.Range("B10").Value = 181
$ws.Range("C10").Value = 344
$ws.Range("D10").Value = 505
$ws.Range("E10").Value = 451
$ws.Range("F10").Value = 325
$ws.Range("G10").Value = 722
$ws.Range("H10").Value = 1003
$ws.Range("I10").Value = 1061
